# Updated cryptos list (price/volume refresh), mirroring the GitHub Actions
# scheduled-job commit. Only the Price (D) and Volume(1h) (E) columns change
# for most rows; rows 16/17 and 24/25 additionally swap rank position
# (coin name/link/price/volume move to the other row).
#
# Note: several Price values are decimal-looking text (e.g. "1.00", "0.165")
# that must stay as literal strings (matching the original inlineStr cells)
# rather than being coerced to numbers by Excel - those are written with a
# leading apostrophe to force text entry, same as typing them in the UI.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "59.556.14"
$ws.Range("E2").Value = "  +0.60%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.606.79"
$ws.Range("E3").Value = "  +0.32%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.02%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'537.73"
$ws.Range("E5").Value = "  +2.91%  "

# Row 6 - Solana
$ws.Range("D6").Value = "'141.55"
$ws.Range("E6").Value = "  +1.63%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.02%  "

# Row 8 - XRP
$ws.Range("E8").Value = "  +0.15%  "

# Row 9 - Toncoin
$ws.Range("E9").Value = "  -0.45%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  +1.14%  "

# Row 11 - Cardano
$ws.Range("E11").Value = "  +1.47%  "

# Row 12 - TRON
$ws.Range("E12").Value = "  -1.06%  "

# Row 13 - WrappedliquidstakedEther2.0
$ws.Range("D13").Value = "3.061.18"
$ws.Range("E13").Value = "  +0.36%  "

# Row 14 - WrappedBTC
$ws.Range("D14").Value = "59.457.44"
$ws.Range("E14").Value = "  +0.59%  "

# Row 15 - Avalanche
$ws.Range("E15").Value = "  +1.19%  "

# Row 16 - was ShibaInu, now WrappedEther (swapped with row 17)
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "2.649.18"
$ws.Range("E16").Value = "  +2.37%  "

# Row 17 - was WrappedEther, now ShibaInu (swapped with row 16)
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").Value = "'0.0000133"
$ws.Range("E17").Value = "  +0.52%  "

# Row 18 - BitcoinCash
$ws.Range("D18").Value = "'340.66"
$ws.Range("E18").Value = "  -0.13%  "

# Row 20 - Chainlink
$ws.Range("E20").Value = "  +0.18%  "

# Row 21 - Uniswap
$ws.Range("E21").Value = "  -2.16%  "

# Row 22 - Dai
$ws.Range("D22").Value = "'1.00"

# Row 23 - Litecoin
$ws.Range("D23").Value = "'67.18"
$ws.Range("E23").Value = "  +1.03%  "

# Row 24 - was Kaspa, now Polygon (swapped with row 25)
$ws.Range("B24").Value = "Polygon"
$ws.Range("C24").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D24").Value = "'0.408"
$ws.Range("E24").Value = "  +0.91%  "

# Row 25 - was Polygon, now Kaspa (swapped with row 24)
$ws.Range("B25").Value = "Kaspa"
$ws.Range("C25").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D25").Value = "'0.165"
$ws.Range("E25").Value = "  -1.35%  "

# Row 26 - Binance-PegBSC-USD
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  +0.21%  "

# Row 27 - InternetComputer(DFINITY)
$ws.Range("E27").Value = "  +2.18%  "

# Row 28 - PEPE
$ws.Range("D28").Value = "0.0₃0742"

# Row 29 - USDe
$ws.Range("E29").Value = "  +0.01%  "

# Row 30 - PancakeSwap
$ws.Range("E30").Value = "  +5.46%  "

# Row 31 - Aptos
$ws.Range("E31").Value = "  -0.73%  "

# Row 32 - EthereumClassic
$ws.Range("E32").Value = "  +0.26%  "

# Row 33 - Monero
$ws.Range("D33").Value = "'150.42"
$ws.Range("E33").Value = "  +0.82%  "

# Row 34 - NEARProtocol
$ws.Range("E34").Value = "  +0.52%  "

# Row 35 - ImmutableX
$ws.Range("E35").Value = "  +0.39%  "

# Row 36 - SuiNetwork
$ws.Range("D36").Value = "'0.843"
$ws.Range("E36").Value = "  +3.68%  "

# Row 37 - Stacks
$ws.Range("E37").Value = "  -0.99%  "

# Row 38 - Fetch.AI
$ws.Range("E38").Value = "  -0.12%  "

# Row 39 - Filecoin
$ws.Range("E39").Value = "  +0.47%  "

# Row 40 - FirstDigitalUSD
$ws.Range("D40").Value = "'0.999"
$ws.Range("E40").Value = "  -0.06%  "

# Row 41 - Bittensor
$ws.Range("D41").Value = "'274.88"
$ws.Range("E41").Value = "  +1.15%  "

# Row 42 - Mantle
$ws.Range("E42").Value = "  -0.05%  "

# Row 43 - WhiteBITCoin
$ws.Range("D43").Value = "'10.73"
$ws.Range("E43").Value = "  -0.21%  "

# Row 44 - Stellar
$ws.Range("D44").Value = "'0.0948"
$ws.Range("E44").Value = "  -0.34%  "

# Row 45 - Hedera
$ws.Range("E45").Value = "  +1.43%  "

# Row 46 - VeChain
$ws.Range("E46").Value = "  +0.82%  "

# Row 47 - Maker
$ws.Range("D47").Value = "1.936.52"

# Row 48 - InjectiveProtocol
$ws.Range("D48").Value = "'18.40"
$ws.Range("E48").Value = "  +2.26%  "

# Row 49 - RenderToken
$ws.Range("E49").Value = "  +0.82%  "

# Row 50 - Aave
$ws.Range("D50").Value = "'110.87"
$ws.Range("E50").Value = "  -2.64%  "

# Row 51 - ZEEBU
$ws.Range("E51").Value = "  +1.84%  "
